# Update Bahamut_Profits market-price / leve-profit figures
# (values refreshed by the scheduled market-data runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1402.6477
$ws.Range("I15").Value = 1402.6477
$ws.Range("K15").Value = 4207.9431
$ws.Range("M15").Value = -4038.9431
$ws.Range("H28").Value = 1438.6316
$ws.Range("I28").Value = 1981
$ws.Range("J28").Value = 508.85715
$ws.Range("K28").Value = 1981
$ws.Range("L28").Value = 508.85715
$ws.Range("M28").Value = -1496
$ws.Range("N28").Value = -1478.85715
$ws.Range("H33").Value = 38613.152
$ws.Range("I33").Value = 40145.68
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 40145.68
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -39916.68
$ws.Range("N33").Value = -758
$ws.Range("H43").Value = 1251345
$ws.Range("I43").Value = 1894.5
$ws.Range("J43").Value = 2500795.5
$ws.Range("K43").Value = 1894.5
$ws.Range("L43").Value = 2500795.5
$ws.Range("M43").Value = -1825.5
$ws.Range("N43").Value = -2500933.5
$ws.Range("H106").Value = 2862.7856
$ws.Range("I106").Value = 2733.2222
$ws.Range("K106").Value = 2733.2222
$ws.Range("M106").Value = -2102.2222
$ws.Range("H121").Value = 672.96875
$ws.Range("J121").Value = 664.5
$ws.Range("L121").Value = 1993.5
$ws.Range("N121").Value = -5487.5
$ws.Range("H125").Value = 20835696
$ws.Range("I125").Value = 29412558
$ws.Range("J125").Value = 6176.4287
$ws.Range("K125").Value = 264713022
$ws.Range("L125").Value = 55587.85830000001
$ws.Range("M125").Value = -264710562
$ws.Range("N125").Value = -60507.85830000001
$ws.Range("H137").Value = 8773418
$ws.Range("I137").Value = 1447.8055
$ws.Range("J137").Value = 23811082
$ws.Range("K137").Value = 4343.416499999999
$ws.Range("L137").Value = 71433246
$ws.Range("M137").Value = -1793.416499999999
$ws.Range("N137").Value = -71438346

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 50391.81
$ws.Range("I2").Value = 60713.41
$ws.Range("K2").Value = 60713.41
$ws.Range("M2").Value = -60600.41
$ws.Range("H32").Value = 24570.348
$ws.Range("I32").Value = 21119.422
$ws.Range("K32").Value = 21119.422
$ws.Range("M32").Value = -20832.422
$ws.Range("H45").Value = 6277.92
$ws.Range("I45").Value = 7365.8823
$ws.Range("J45").Value = 3966
$ws.Range("K45").Value = 7365.8823
$ws.Range("L45").Value = 3966
$ws.Range("M45").Value = -6988.8823
$ws.Range("N45").Value = -4720
$ws.Range("H63").Value = 1985.8928
$ws.Range("I63").Value = 1985.3704
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 1985.3704
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1299.3704
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 1985.8928
$ws.Range("I66").Value = 1985.3704
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 9926.852000000001
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -6494.852000000001
$ws.Range("N66").Value = -16864
$ws.Range("H110").Value = 510.45456
$ws.Range("I110").Value = 505.19232
$ws.Range("J110").Value = 530
$ws.Range("K110").Value = 505.19232
$ws.Range("L110").Value = 530
$ws.Range("M110").Value = 1539.80768
$ws.Range("N110").Value = -4620
$ws.Range("H116").Value = 50391.81
$ws.Range("I116").Value = 60713.41
$ws.Range("K116").Value = 60713.41
$ws.Range("M116").Value = -58419.41
$ws.Range("H132").Value = 2137.8684
$ws.Range("I132").Value = 1578.9131
$ws.Range("J132").Value = 2994.9333
$ws.Range("K132").Value = 4736.7393
$ws.Range("L132").Value = 8984.7999
$ws.Range("M132").Value = -2206.7393
$ws.Range("N132").Value = -14044.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 50391.81
$ws.Range("I3").Value = 60713.41
$ws.Range("K3").Value = 60713.41
$ws.Range("M3").Value = -60599.41
$ws.Range("H94").Value = 1279.3158
$ws.Range("I94").Value = 1186.3334
$ws.Range("J94").Value = 1628
$ws.Range("K94").Value = 1186.3334
$ws.Range("L94").Value = 1628
$ws.Range("M94").Value = -735.3334
$ws.Range("N94").Value = -2530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1387.5084
$ws.Range("I31").Value = 802.37933
$ws.Range("J31").Value = 1953.1333
$ws.Range("K31").Value = 802.37933
$ws.Range("L31").Value = 1953.1333
$ws.Range("M31").Value = -507.37933
$ws.Range("N31").Value = -2543.1333
$ws.Range("H34").Value = 1387.5084
$ws.Range("I34").Value = 802.37933
$ws.Range("J34").Value = 1953.1333
$ws.Range("K34").Value = 802.37933
$ws.Range("L34").Value = 1953.1333
$ws.Range("M34").Value = -600.37933
$ws.Range("N34").Value = -2357.1333
$ws.Range("H58").Value = 10631
$ws.Range("I58").Value = 1427
$ws.Range("J58").Value = 19835
$ws.Range("K58").Value = 1427
$ws.Range("L58").Value = 19835
$ws.Range("M58").Value = -1224
$ws.Range("N58").Value = -20241
$ws.Range("H86").Value = 3222.7144
$ws.Range("I86").Value = 3205.818
$ws.Range("J86").Value = 3284.6667
$ws.Range("K86").Value = 3205.818
$ws.Range("L86").Value = 3284.6667
$ws.Range("M86").Value = -2082.818
$ws.Range("N86").Value = -5530.6667
$ws.Range("H89").Value = 3222.7144
$ws.Range("I89").Value = 3205.818
$ws.Range("J89").Value = 3284.6667
$ws.Range("K89").Value = 16029.09
$ws.Range("L89").Value = 16423.3335
$ws.Range("M89").Value = -10413.09
$ws.Range("N89").Value = -27655.3335
$ws.Range("H132").Value = 3615.25
$ws.Range("I132").Value = 3022.4
$ws.Range("J132").Value = 4299.3076
$ws.Range("K132").Value = 9067.200000000001
$ws.Range("L132").Value = 12897.9228
$ws.Range("M132").Value = -6537.200000000001
$ws.Range("N132").Value = -17957.9228
$ws.Range("H134").Value = 2627.5715
$ws.Range("I134").Value = 1931.7778
$ws.Range("K134").Value = 5795.3334
$ws.Range("M134").Value = -3260.3334
$ws.Range("H136").Value = 10631
$ws.Range("I136").Value = 1427
$ws.Range("J136").Value = 19835
$ws.Range("K136").Value = 4281
$ws.Range("L136").Value = 59505
$ws.Range("M136").Value = -1731
$ws.Range("N136").Value = -64605
$ws.Range("H140").Value = 47948
$ws.Range("J140").Value = 47948
$ws.Range("L140").Value = 47948
$ws.Range("N140").Value = -58308

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 90.666664
$ws.Range("I40").Value = 93.25
$ws.Range("J40").Value = 70
$ws.Range("K40").Value = 373
$ws.Range("L40").Value = 280
$ws.Range("M40").Value = -304
$ws.Range("N40").Value = -418
$ws.Range("H68").Value = 1370.0934
$ws.Range("I68").Value = 744.26666
$ws.Range("K68").Value = 2232.79998
$ws.Range("M68").Value = -1421.79998
$ws.Range("H71").Value = 1370.0934
$ws.Range("I71").Value = 744.26666
$ws.Range("K71").Value = 6698.39994
$ws.Range("M71").Value = -2642.39994
$ws.Range("H107").Value = 836
$ws.Range("J107").Value = 1542.2693
$ws.Range("L107").Value = 4626.8079
$ws.Range("N107").Value = -8466.8079
$ws.Range("H131").Value = 17657.127
$ws.Range("I131").Value = 84687.5
$ws.Range("J131").Value = 1885.2745
$ws.Range("K131").Value = 254062.5
$ws.Range("L131").Value = 5655.8235
$ws.Range("M131").Value = -249022.5
$ws.Range("N131").Value = -15735.8235
$ws.Range("H134").Value = 7532.636

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3400
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 3400
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 3400
$ws.Range("N102").Value = -6644
$ws.Range("M102").ClearContents()
$ws.Range("H113").Value = 1479.3636
$ws.Range("I113").Value = 1157.1428
$ws.Range("K113").Value = 1157.1428
$ws.Range("M113").Value = 1012.8572
$ws.Range("H132").Value = 3065.4783
$ws.Range("I132").Value = 2166.6667
$ws.Range("J132").Value = 4046
$ws.Range("K132").Value = 6500.000100000001
$ws.Range("L132").Value = 12138
$ws.Range("M132").Value = -3970.000100000001
$ws.Range("N132").Value = -17198
$ws.Range("H138").Value = 25009.908
$ws.Range("J138").Value = 25009.908
$ws.Range("L138").Value = 25009.908
$ws.Range("N138").Value = -35289.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 25125.273
$ws.Range("I92").Value = 24300
$ws.Range("K92").Value = 24300
$ws.Range("M92").Value = -21804
$ws.Range("H132").Value = 2035358.9
$ws.Range("I132").Value = 2876252.2
$ws.Range("J132").Value = 3199.75
$ws.Range("K132").Value = 8628756.600000001
$ws.Range("L132").Value = 9599.25
$ws.Range("M132").Value = -8626226.600000001
$ws.Range("N132").Value = -14659.25
$ws.Range("H139").Value = 48402.5
$ws.Range("J139").Value = 48402.5
$ws.Range("L139").Value = 48402.5
$ws.Range("N139").Value = -58682.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 502.4
$ws.Range("I113").Value = 532.4286
$ws.Range("J113").Value = 432.33334
$ws.Range("K113").Value = 1597.2858
$ws.Range("L113").Value = 1297.00002
$ws.Range("M113").Value = 572.7142000000001
$ws.Range("N113").Value = -5637.000019999999
$ws.Range("H136").Value = 2944.1462
$ws.Range("I136").Value = 3860.7144
$ws.Range("J136").Value = 1981.75
$ws.Range("K136").Value = 11582.1432
$ws.Range("L136").Value = 5945.25
$ws.Range("M136").Value = -9032.143199999999
$ws.Range("N136").Value = -11045.25
$ws.Range("H138").Value = 37621.8
$ws.Range("J138").Value = 37621.8
$ws.Range("L138").Value = 37621.8
$ws.Range("N138").Value = -47901.8
